$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.057.34'
$ws.Range("E2").Value = '  +4.29%  '
$ws.Range("D3").Value = '3.532.56'
$ws.Range("E3").Value = '  +14.95%  '
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("D5").Value = '586.93'
$ws.Range("E5").Value = '  +2.30%  '
$ws.Range("D6").Value = '183.43'
$ws.Range("E6").Value = '  +8.31%  '
$ws.Range("D7").Value = '3.524.59'
$ws.Range("E7").Value = '  +14.78%  '
$ws.Range("D8").Value = '1.00'
$ws.Range("E8").Value = '  -0.01%  '
$ws.Range("E9").Value = '  +3.89%  '
$ws.Range("D10").Value = '6.54'
$ws.Range("E10").Value = '  +3.97%  '
$ws.Range("E11").Value = '  +5.59%  '
$ws.Range("D12").Value = '0.488'
$ws.Range("E12").Value = '  +4.28%  '
$ws.Range("E13").Value = '  +4.47%  '
$ws.Range("D14").Value = '38.17'
$ws.Range("E14").Value = '  +6.83%  '
$ws.Range("D15").Value = '4.112.75'
$ws.Range("E15").Value = '  +14.65%  '
$ws.Range("D16").Value = '69.325.83'
$ws.Range("E16").Value = '  +4.78%  '
$ws.Range("E17").Value = '  +1.34%  '
$ws.Range("D18").Value = '3.508.85'
$ws.Range("E18").Value = '  +14.11%  '
$ws.Range("D19").Value = '7.40'
$ws.Range("E19").Value = '  +6.81%  '
$ws.Range("D20").Value = '16.74'
$ws.Range("E20").Value = '  +2.04%  '
$ws.Range("D21").Value = '502.13'
$ws.Range("E21").Value = '  +3.80%  '
$ws.Range("D22").Value = '9.04'
$ws.Range("E22").Value = '  +18.42%  '
$ws.Range("D23").Value = '0.734'
$ws.Range("E23").Value = '  +7.31%  '
$ws.Range("D24").Value = '86.07'
$ws.Range("E24").Value = '  +4.52%  '
$ws.Range("D25").Value = '13.28'
$ws.Range("E25").Value = '  +5.19%  '
$ws.Range("D26").Value = '2.35'
$ws.Range("E26").Value = '  +6.85%  '
$ws.Range("D27").Value = '10.57'
$ws.Range("E27").Value = '  +4.26%  '
$ws.Range("E29").Value = '  +11.68%  '
$ws.Range("D30").Value = '8.03'
$ws.Range("E30").Value = '  +2.44%  '
$ws.Range("D31").Value = '31.27'
$ws.Range("E31").Value = '  +13.00%  '
$ws.Range("B32").Value = 'PancakeSwap'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D32").Value = '2.71'
$ws.Range("E32").Value = '  +4.65%  '
$ws.Range("B33").Value = 'PEPE'
$ws.Range("C33").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D33").Value = '0.0000108'
$ws.Range("E33").Value = '  +20.06%  '
$ws.Range("D34").Value = '0.116'
$ws.Range("E34").Value = '  +4.90%  '
$ws.Range("D35").Value = '0.999'
$ws.Range("E35").Value = '  -0.03%  '
$ws.Range("E36").Value = '  +9.08%  '
$ws.Range("E37").Value = '  +7.60%  '
$ws.Range("E38").Value = '  +10.51%  '
$ws.Range("D39").Value = '2.07'
$ws.Range("E39").Value = '  +6.44%  '
$ws.Range("D40").Value = '46.32'
$ws.Range("E40").Value = '  -2.22%  '
$ws.Range("D41").Value = '50.53'
$ws.Range("E41").Value = '  +3.10%  '
$ws.Range("E42").Value = '  +3.37%  '
$ws.Range("D43").Value = '8.69'
$ws.Range("E43").Value = '  +5.97%  '
$ws.Range("D44").Value = '3.007.42'
$ws.Range("E44").Value = '  +8.51%  '
$ws.Range("D45").Value = '2.79'
$ws.Range("E45").Value = '  +11.29%  '
$ws.Range("D46").Value = '395.02'
$ws.Range("E46").Value = '  +9.02%  '
$ws.Range("E47").Value = '  +5.30%  '
$ws.Range("E48").Value = '  +13.82%  '
$ws.Range("D49").Value = '134.16'
$ws.Range("E49").Value = '  -0.35%  '
$ws.Range("D51").Value = '2.45'
$ws.Range("E51").Value = '  +14.53%  '
